$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44550   # D2 Fecha
$ws.Cells.Item(2, 10).Value = 1300  # J2 Volumen
$ws.Cells.Item(2, 11).Value = 1000  # K2 Precio minimo
$ws.Cells.Item(2, 12).Value = 1200  # L2 Precio maximo
$ws.Cells.Item(2, 13).Value = 1100  # M2 Precio promedio ponderado
$ws.Cells.Item(2, 16).Value = 1100  # P2 Precio $/Kg

$ws.Cells.Item(3, 4).Value = 44883   # D3 Fecha
$ws.Cells.Item(3, 10).Value = 800  # J3 Volumen
$ws.Cells.Item(3, 11).Value = 550  # K3 Precio minimo
$ws.Cells.Item(3, 12).Value = 600  # L3 Precio maximo
$ws.Cells.Item(3, 13).Value = 575  # M3 Precio promedio ponderado
$ws.Cells.Item(3, 16).Value = 575  # P3 Precio $/Kg

$ws.Cells.Item(4, 4).Value = 44453   # D4 Fecha
$ws.Cells.Item(4, 10).Value = 1000  # J4 Volumen
$ws.Cells.Item(4, 11).Value = 800  # K4 Precio minimo
$ws.Cells.Item(4, 12).Value = 900  # L4 Precio maximo
$ws.Cells.Item(4, 13).Value = 850  # M4 Precio promedio ponderado
$ws.Cells.Item(4, 16).Value = 850  # P4 Precio $/Kg

$ws.Cells.Item(5, 4).Value = 44175   # D5 Fecha
$ws.Cells.Item(5, 10).Value = 1600  # J5 Volumen
$ws.Cells.Item(5, 11).Value = 1000  # K5 Precio minimo
$ws.Cells.Item(5, 12).Value = 1200  # L5 Precio maximo
$ws.Cells.Item(5, 13).Value = 1100  # M5 Precio promedio ponderado
$ws.Cells.Item(5, 16).Value = 1100  # P5 Precio $/Kg

$ws.Cells.Item(6, 4).Value = 44649   # D6 Fecha
$ws.Cells.Item(6, 10).Value = 600  # J6 Volumen
$ws.Cells.Item(6, 11).Value = 900  # K6 Precio minimo
$ws.Cells.Item(6, 12).Value = 1000  # L6 Precio maximo
$ws.Cells.Item(6, 13).Value = 950  # M6 Precio promedio ponderado
$ws.Cells.Item(6, 16).Value = 950  # P6 Precio $/Kg

$ws.Cells.Item(7, 4).Value = 44687   # D7 Fecha
$ws.Cells.Item(7, 10).Value = 1000  # J7 Volumen
$ws.Cells.Item(7, 11).Value = 1200  # K7 Precio minimo
$ws.Cells.Item(7, 12).Value = 1300  # L7 Precio maximo
$ws.Cells.Item(7, 13).Value = 1250  # M7 Precio promedio ponderado
$ws.Cells.Item(7, 16).Value = 1250  # P7 Precio $/Kg

$ws.Cells.Item(8, 4).Value = 44442   # D8 Fecha
$ws.Cells.Item(8, 10).Value = 1250  # J8 Volumen
$ws.Cells.Item(8, 11).Value = 850  # K8 Precio minimo
$ws.Cells.Item(8, 12).Value = 900  # L8 Precio maximo
$ws.Cells.Item(8, 13).Value = 875  # M8 Precio promedio ponderado
$ws.Cells.Item(8, 16).Value = 875  # P8 Precio $/Kg

$ws.Cells.Item(9, 4).Value = 44784   # D9 Fecha
$ws.Cells.Item(9, 10).Value = 1000  # J9 Volumen
$ws.Cells.Item(9, 11).Value = 1200  # K9 Precio minimo
$ws.Cells.Item(9, 12).Value = 1300  # L9 Precio maximo
$ws.Cells.Item(9, 13).Value = 1250  # M9 Precio promedio ponderado
$ws.Cells.Item(9, 16).Value = 1250  # P9 Precio $/Kg

$ws.Cells.Item(10, 4).Value = 44476   # D10 Fecha
$ws.Cells.Item(10, 10).Value = 900  # J10 Volumen
$ws.Cells.Item(10, 11).Value = 700  # K10 Precio minimo
$ws.Cells.Item(10, 12).Value = 800  # L10 Precio maximo
$ws.Cells.Item(10, 13).Value = 750  # M10 Precio promedio ponderado
$ws.Cells.Item(10, 16).Value = 750  # P10 Precio $/Kg

$ws.Cells.Item(11, 4).Value = 44407   # D11 Fecha
$ws.Cells.Item(11, 10).Value = 1000  # J11 Volumen
$ws.Cells.Item(11, 11).Value = 1200  # K11 Precio minimo
$ws.Cells.Item(11, 12).Value = 1300  # L11 Precio maximo
$ws.Cells.Item(11, 13).Value = 1250  # M11 Precio promedio ponderado
$ws.Cells.Item(11, 16).Value = 1250  # P11 Precio $/Kg

$ws.Cells.Item(12, 4).Value = 44243   # D12 Fecha
$ws.Cells.Item(12, 10).Value = 1200  # J12 Volumen
$ws.Cells.Item(12, 11).Value = 1200  # K12 Precio minimo
$ws.Cells.Item(12, 12).Value = 1300  # L12 Precio maximo
$ws.Cells.Item(12, 13).Value = 1250  # M12 Precio promedio ponderado
$ws.Cells.Item(12, 16).Value = 1250  # P12 Precio $/Kg

$ws.Cells.Item(13, 4).Value = 44638   # D13 Fecha
$ws.Cells.Item(13, 10).Value = 1000  # J13 Volumen
$ws.Cells.Item(13, 11).Value = 900  # K13 Precio minimo
$ws.Cells.Item(13, 12).Value = 950  # L13 Precio maximo
$ws.Cells.Item(13, 13).Value = 925  # M13 Precio promedio ponderado
$ws.Cells.Item(13, 16).Value = 925  # P13 Precio $/Kg

$ws.Cells.Item(14, 4).Value = 44449   # D14 Fecha
$ws.Cells.Item(14, 10).Value = 1300  # J14 Volumen
$ws.Cells.Item(14, 11).Value = 900  # K14 Precio minimo
$ws.Cells.Item(14, 12).Value = 950  # L14 Precio maximo
$ws.Cells.Item(14, 13).Value = 925  # M14 Precio promedio ponderado
$ws.Cells.Item(14, 16).Value = 925  # P14 Precio $/Kg

$ws.Cells.Item(15, 4).Value = 44229   # D15 Fecha
$ws.Cells.Item(15, 10).Value = 1500  # J15 Volumen
$ws.Cells.Item(15, 11).Value = 1400  # K15 Precio minimo
$ws.Cells.Item(15, 12).Value = 1500  # L15 Precio maximo
$ws.Cells.Item(15, 13).Value = 1450  # M15 Precio promedio ponderado
$ws.Cells.Item(15, 16).Value = 1450  # P15 Precio $/Kg

$ws.Cells.Item(16, 4).Value = 44291   # D16 Fecha
$ws.Cells.Item(16, 10).Value = 1000  # J16 Volumen
$ws.Cells.Item(16, 11).Value = 1000  # K16 Precio minimo
$ws.Cells.Item(16, 12).Value = 1200  # L16 Precio maximo
$ws.Cells.Item(16, 13).Value = 1100  # M16 Precio promedio ponderado
$ws.Cells.Item(16, 16).Value = 1100  # P16 Precio $/Kg

$ws.Cells.Item(17, 4).Value = 44284   # D17 Fecha
$ws.Cells.Item(17, 10).Value = 1500  # J17 Volumen
$ws.Cells.Item(17, 11).Value = 800  # K17 Precio minimo
$ws.Cells.Item(17, 12).Value = 850  # L17 Precio maximo
$ws.Cells.Item(17, 13).Value = 825  # M17 Precio promedio ponderado
$ws.Cells.Item(17, 16).Value = 825  # P17 Precio $/Kg

$ws.Cells.Item(18, 4).Value = 44607   # D18 Fecha
$ws.Cells.Item(18, 10).Value = 900  # J18 Volumen
$ws.Cells.Item(18, 11).Value = 1300  # K18 Precio minimo
$ws.Cells.Item(18, 12).Value = 1400  # L18 Precio maximo
$ws.Cells.Item(18, 13).Value = 1350  # M18 Precio promedio ponderado
$ws.Cells.Item(18, 16).Value = 1350  # P18 Precio $/Kg

$ws.Cells.Item(19, 4).Value = 44673   # D19 Fecha
$ws.Cells.Item(19, 10).Value = 900  # J19 Volumen
$ws.Cells.Item(19, 11).Value = 1300  # K19 Precio minimo
$ws.Cells.Item(19, 12).Value = 1400  # L19 Precio maximo
$ws.Cells.Item(19, 13).Value = 1350  # M19 Precio promedio ponderado
$ws.Cells.Item(19, 16).Value = 1350  # P19 Precio $/Kg

$ws.Cells.Item(20, 4).Value = 44656   # D20 Fecha
$ws.Cells.Item(20, 10).Value = 1000  # J20 Volumen
$ws.Cells.Item(20, 11).Value = 900  # K20 Precio minimo
$ws.Cells.Item(20, 12).Value = 1000  # L20 Precio maximo
$ws.Cells.Item(20, 13).Value = 950  # M20 Precio promedio ponderado
$ws.Cells.Item(20, 16).Value = 950  # P20 Precio $/Kg

$ws.Cells.Item(21, 4).Value = 44455   # D21 Fecha
$ws.Cells.Item(21, 10).Value = 1100  # J21 Volumen
$ws.Cells.Item(21, 11).Value = 900  # K21 Precio minimo
$ws.Cells.Item(21, 12).Value = 1000  # L21 Precio maximo
$ws.Cells.Item(21, 13).Value = 950  # M21 Precio promedio ponderado
$ws.Cells.Item(21, 16).Value = 950  # P21 Precio $/Kg

$ws.Cells.Item(22, 4).Value = 44484   # D22 Fecha
$ws.Cells.Item(22, 10).Value = 900  # J22 Volumen
$ws.Cells.Item(22, 11).Value = 750  # K22 Precio minimo
$ws.Cells.Item(22, 12).Value = 800  # L22 Precio maximo
$ws.Cells.Item(22, 13).Value = 775  # M22 Precio promedio ponderado
$ws.Cells.Item(22, 16).Value = 775  # P22 Precio $/Kg

$ws.Cells.Item(23, 4).Value = 44341   # D23 Fecha
$ws.Cells.Item(23, 10).Value = 1300  # J23 Volumen
$ws.Cells.Item(23, 11).Value = 900  # K23 Precio minimo
$ws.Cells.Item(23, 12).Value = 1000  # L23 Precio maximo
$ws.Cells.Item(23, 13).Value = 950  # M23 Precio promedio ponderado
$ws.Cells.Item(23, 16).Value = 950  # P23 Precio $/Kg
